$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 13 data
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "search for component"
$ws.Range("C13").Value = "verify image in component table only searched once"

# Freeze panes at row 2 (split after header row)
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Select C13 as the active cell after freezing
[void]$ws.Range("C13").Select()

# Adjust column C width to fit new content (closest reachable value to the
# recorded best-fit width of 42.05078125 given COM's pixel-grid rounding)
$ws.Columns.Item(3).ColumnWidth = 41.1
